# Insert a new weekly record for "Vega Monumental Concepción - Berenjena" at row 14,
# pushing the existing rows 14-71 down to rows 15-72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 44607
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100112001
$ws.Cells.Item(14, 7).Value = "Berenjena"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 100
$ws.Cells.Item(14, 11).Value = 13000
$ws.Cells.Item(14, 12).Value = 14000
$ws.Cells.Item(14, 13).Value = 13500
$ws.Cells.Item(14, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(14, 15).Value = "Región Metropolitana"
$ws.Cells.Item(14, 16).Value = 225
$ws.Cells.Item(14, 17).Value = 60
$ws.Cells.Item(14, 18).Value = "Hortaliza"

$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
